$d = $word.ActiveDocument

# --- Edit 1: "N=15, p=0.2 = >1-p = 0.8" -> "N=16, p=0.19 = >1-p = 0.81" ---
# Scope the Find to the specific paragraph so we don't touch similar numbers
# ("15", "2") appearing elsewhere in the document (e.g. subnet examples).
$paras = $d.Paragraphs
$targetPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith("N=15, p=0.2")) {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $rng = $targetPara.Range
    $rng.Find.Execute("N=15, p=0.2 = >1-p = 0.8", $true, $true, $false, $false, $false, $true, 1, $false, "N=16, p=0.19 = >1-p = 0.81", 2)
}

# --- Edit 2: first "B" answer (for 192.24.6.0) -> "D" ---
# There are three single-letter answer paragraphs ("B", "B", "D"); only the
# first one (right after the 192.24.6.0 line) changes.
$paras2 = $d.Paragraphs
$firstB = $null
for ($i = 1; $i -le $paras2.Count; $i++) {
    $p = $paras2.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "B") {
        $firstB = $p
        break
    }
}

if ($firstB -ne $null) {
    $rng2 = $firstB.Range
    $rng2.Find.Execute("B", $true, $true, $false, $false, $false, $true, 1, $false, "D", 2)
}
